$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the measurement values that used to live in B2:F4, leaving the
# "A" label cells in column A untouched (mirrors the data having been wiped
# out while the row/label structure stays intact).
$ws.Range("B2:F4").ClearContents()

# Reflect the new selection left behind after clearing that range
# (activeCell="B2", sqref="B2:F4").
$ws.Range("B2:F4").Select() | Out-Null
